# Refactoring 9/28/24 @ 19:07
#
# 1) Rename sheets: RequestSignup -> Account, RequestLogin -> Auth
# 2) Update Auth (sheet2) sample row: username/testuser40 -> email/testuser40@gmail.com
# 3) Add four new sheets (Email, Password, Admin, Change) with header + sample rows,
#    reusing existing styles/shared-strings wherever the source workbook already has
#    a matching (style, text) cell so the style table & shared-string table line up
#    with the target.

$wb = $excel.ActiveWorkbook

$account = $wb.Worksheets.Item(1)   # was "RequestSignup"
$auth    = $wb.Worksheets.Item(2)   # was "RequestLogin"
$data    = $wb.Worksheets.Item(3)   # "Data" (unchanged)
$session = $wb.Worksheets.Item(4)   # "Session" (unchanged)

# ---------------------------------------------------------------------------
# 1) Sheet renames
# ---------------------------------------------------------------------------
$account.Name = "Account"
$auth.Name = "Auth"

# ---------------------------------------------------------------------------
# 2) Auth sheet sample-row update: A1 username->email, A2 testuser40->testuser40@gmail.com
#    (style stays exactly as-is -- only the text changes -- so just overwrite .Value)
# ---------------------------------------------------------------------------
$auth.Range("A1").Value = "email"
$auth.Range("A2").Value = "testuser40@gmail.com"

# ---------------------------------------------------------------------------
# 3) Create the four new sheets, positioned after "Session"
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$email = $wb.Worksheets.Add($null, $last)
$email.Name = "Email"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$password = $wb.Worksheets.Add($null, $last)
$password.Name = "Password"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$admin = $wb.Worksheets.Add($null, $last)
$admin.Name = "Admin"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$change = $wb.Worksheets.Add($null, $last)
$change.Name = "Change"

# ---------------------------------------------------------------------------
# Email sheet: email | password | newEmail | removeEmail
# ---------------------------------------------------------------------------
$data.Range("B1").Copy($email.Range("A1"))          # email
$auth.Range("B1").Copy($email.Range("B1"))           # password
$data.Range("B1").Copy($email.Range("C1"))          # newEmail (copy style, fix text below)
$email.Range("C1").Value = "newEmail"
$data.Range("B1").Copy($email.Range("D1"))          # removeEmail (copy style, fix text below)
$email.Range("D1").Value = "removeEmail"

$auth.Range("B2").Copy($email.Range("A2"))           # testuser40@gmail.com
$auth.Range("C2").Copy($email.Range("B2"))           # Test1234$
$auth.Range("B2").Copy($email.Range("C2"))           # testuser41@gmail.com (copy style, fix text below)
$email.Range("C2").Value = "testuser41@gmail.com"
$auth.Range("B2").Copy($email.Range("D2"))           # testuser41@gmail.com (copy style, fix text below)
$email.Range("D2").Value = "testuser41@gmail.com"

# ---------------------------------------------------------------------------
# Password sheet: email | username | oldPassword | password | confirmPassword |
#                  confirmationCode | accessToken | refreshToken
# ---------------------------------------------------------------------------
$data.Range("B1").Copy($password.Range("A1"))       # email
$data.Range("C1").Copy($password.Range("B1"))       # username
$data.Range("B1").Copy($password.Range("C1"))       # oldPassword (copy style, fix text below)
$password.Range("C1").Value = "oldPassword"
$auth.Range("B1").Copy($password.Range("D1"))        # password
$account.Range("D1").Copy($password.Range("E1"))        # confirmPassword
$data.Range("I1").Copy($password.Range("F1"))       # confirmationCode
$session.Range("U1").Copy($password.Range("G1"))     # accessToken
$session.Range("V1").Copy($password.Range("H1"))     # refreshToken

# A2 uses a brand-new style (font "docs-JetBrains Mono", left aligned) -- create it once
# then propagate via copy so the style table gains exactly one new font + one new xf.
$password.Range("A2").Value = "testuser40@gmail.com"
$password.Range("A2").Font.Name = "docs-JetBrains Mono"
$password.Range("A2").Font.Color = 0
$password.Range("A2").HorizontalAlignment = -4131

$account.Range("B2").Copy($password.Range("B2"))        # testuser40@gmail.comtestuser40 (copy style, fix text below)
$password.Range("B2").Value = "testuser40@gmail.comtestuser40"
$auth.Range("C2").Copy($password.Range("C2"))        # Test1234$
$auth.Range("C2").Copy($password.Range("D2"))        # Test1234$
$auth.Range("C2").Copy($password.Range("E2"))        # Test1234$

$data.Range("B1").Copy($password.Range("F2"))       # 6a70b6e8-... (copy style, fix text below)
$password.Range("F2").Value = "6a70b6e8-2b18-4655-9988-0b5dc3b7d881"
$data.Range("B1").Copy($password.Range("G2"))       # new accessToken JWT (copy style, fix text below)
$password.Range("G2").Value = "eyJhbGciOiJIUzI1NiJ9.eyJsb2NhdGlvbiI6eyJsb25naXR1ZGUiOjEyMzQ1NjcuMCwibGF0aXR1ZGUiOjEyMzQ1NjcuMH0sImlzcyI6Im9yYi1nYXRld2F5Iiwic3ViIjoidGVzdHVzZXI0MCIsImlhdCI6MTcyNzQ0NTY1NiwiZXhwIjoxNzI3NTMyMDU2fQ.XNxak4X2Bt19WFZTu8gM01JrDCaVWhxBVa5m0DHaJNU"
$data.Range("B1").Copy($password.Range("H2"))       # new refreshToken JWT (copy style, fix text below)
$password.Range("H2").Value = "eyJhbGciOiJIUzI1NiJ9.eyJpc3MiOiJvcmItZ2F0ZXdheSIsInN1YiI6InRlc3R1c2VyNDAiLCJpYXQiOjE3Mjc0NDU2NTYsImV4cCI6MTcyNzUzMjA1Nn0.AmYkHEa_cYZsH2KNsc7hIcvuhcdWSnjRZJa3egpvSrU"

# B3 is blank but carries the same new style as A2 -- copy it straight from A2.
$password.Range("A2").Copy($password.Range("B3"))

# Column widths for E:F (approximate -- Excel quantises ColumnWidth to whole pixels)
$password.Columns.Item(5).ColumnWidth = 12.92
$password.Columns.Item(6).ColumnWidth = 13.1

# ---------------------------------------------------------------------------
# Admin sheet: email / testuser40@gmail.com
# ---------------------------------------------------------------------------
$data.Range("B1").Copy($admin.Range("A1"))          # email
$data.Range("B2").Copy($admin.Range("A2"))          # testuser40@gmail.com

# ---------------------------------------------------------------------------
# Change sheet: email | password | username | firstName | lastName | phone |
#                accessToken | refreshToken
# ---------------------------------------------------------------------------
$data.Range("B1").Copy($change.Range("A1"))         # email
$auth.Range("B1").Copy($change.Range("B1"))          # password
$data.Range("C1").Copy($change.Range("C1"))         # username
$data.Range("E1").Copy($change.Range("D1"))         # firstName
$data.Range("F1").Copy($change.Range("E1"))         # lastName
$account.Range("I1").Copy($change.Range("F1"))          # phone
$session.Range("U1").Copy($change.Range("G1"))       # accessToken
$session.Range("V1").Copy($change.Range("H1"))       # refreshToken

$auth.Range("B2").Copy($change.Range("A2"))          # testuser40@gmail.com
$auth.Range("C2").Copy($change.Range("B2"))          # Test1234$
$account.Range("C2").Copy($change.Range("C2"))          # testuser40 (copy style, fix text below)
$change.Range("C2").Value = "testuser40"
$account.Range("G2").Copy($change.Range("D2"))          # Test
$account.Range("H2").Copy($change.Range("E2"))          # One
$account.Range("K2").Copy($change.Range("F2"))          # 14048205065 (copy style, fix value below)
$change.Range("F2").Value = 14048205065

$password.Range("G2").Copy($change.Range("G2"))      # accessToken JWT
$password.Range("H2").Copy($change.Range("H2"))      # refreshToken JWT
